$d = $word.ActiveDocument

function Get-ParagraphByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $text) {
            return $p
        }
    }
    return $null
}

# 1. Remove the _GoBack bookmark from its original location (first paragraph).
#    It will be re-created later at the new edit location, matching how Word
#    auto-tracks the last edit position as you type.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 2. Colour the "CONCLUSION GENERALE ET PERSPECTIVES" heading red (FF0000),
#    on both the run and the paragraph mark (use the whole paragraph Range so
#    the mark's rPr gets the color too, matching Word's "select paragraph,
#    apply font color" behaviour).
$concl = Get-ParagraphByText $d "CONCLUSION GENERALE ET PERSPECTIVES"
$concl.Range.Font.Color = 255

# 3. Split the "Bibliographies et webographies" paragraph into two paragraphs:
#    "Bibliographies " and " webographies", dropping the "et" in between
#    (as if the user placed the cursor after "Bibliographies " and pressed
#    Enter, then deleted "et").
$biblio = Get-ParagraphByText $d "Bibliographies et webographies"
$splitPoint = $biblio.Range.Start + "Bibliographies ".Length
$splitRange = $d.Range($splitPoint, $splitPoint)
$splitRange.InsertParagraphAfter()

# The newly created paragraph (right after the "Bibliographies " paragraph)
# now starts with "et webographies"; remove the leading "et".
$newPara = $biblio.Next()
$etRange = $d.Range($newPara.Range.Start, $newPara.Range.Start + 2)
$etRange.Delete()

# 4. Re-insert the _GoBack bookmark (collapsed) at the start of the new
#    paragraph, i.e. where the edit actually happened.
$bmRange = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)
